$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first row and a new first column, shifting the existing
# grid down-and-right by one row/column (same as the diff shows).
$ws.Rows.Item(1).Insert()
$ws.Columns.Item(1).Insert()

# New header row 1 (B1:J1) -- sequential numbers 1..9 with the style
# previously only used by column K ("style 10", centered horizontally).
$headerValues = @(1,2,3,4,5,6,7,8,9)
for ($i = 0; $i -lt $headerValues.Length; $i++) {
    $col = 2 + $i
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headerValues[$i]
    $cell.HorizontalAlignment = -4108
}

# New first column (A2:A10) -- sequential numbers 1..9.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Clear out the old "x" markers that used to flag individual board
# cells. A few of them (H5, F6, G7) get overwritten with real numbers
# below anyway; the rest (J2, C3, F3) stay blank.
$ws.Range("J2").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("F3").Value = ""

# Fill in the spiral of numbers across D4:H8.
$ws.Range("D4").Value = 17
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 15
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 13

$ws.Range("D5").Value = 18
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 12

$ws.Range("D6").Value = 19
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 11

$ws.Range("D7").Value = 20
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 8
$ws.Range("G7").Value = 9
$ws.Range("H7").Value = 10

$ws.Range("D8").Value = 21
$ws.Range("E8").Value = 22
$ws.Range("F8").Value = 23
$ws.Range("G8").Value = 24
$ws.Range("H8").Value = 25

# Update the selected cell to match the saved view state.
$ws.Range("F6").Select()
